$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for "Piña" (Feria Lagunitas de Puerto
# Montt) dated 44498. It becomes the new row 49, pushing every subsequent
# record down by one row (old row 151 -> new row 152).
$ws.Rows(49).Insert()

$ws.Cells.Item(49, 1).Value = 4
$ws.Cells.Item(49, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(49, 3).Value = "Los Lagos"
$ws.Cells.Item(49, 4).Value = 44498
$ws.Cells.Item(49, 5).Value = 10
$ws.Cells.Item(49, 6).Value = "Fruta"
$ws.Cells.Item(49, 7).Value = 100108
$ws.Cells.Item(49, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(49, 9).Value = 100108005
$ws.Cells.Item(49, 10).Value = "Piña"
$ws.Cells.Item(49, 11).Value = "Caramelo"
$ws.Cells.Item(49, 12).Value = "Segunda"
$ws.Cells.Item(49, 13).Value = 360
$ws.Cells.Item(49, 14).Value = 20000
$ws.Cells.Item(49, 15).Value = 21000
$ws.Cells.Item(49, 16).Value = 20500
$ws.Cells.Item(49, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(49, 18).Value = "Ecuador"
$ws.Cells.Item(49, 19).Value = 1464
$ws.Cells.Item(49, 20).Value = 14
